# Auto-generated update of scheduled market-price data (currentAveragePrice /
# NQ/HQ prices / leve profit calculations) across all craft sheets.
$wb = $excel.ActiveWorkbook

$sheetUpdates = @{}

$sheetUpdates["ALC"] = @{
    "H6" = 271.1389
    "I6" = 147.6
    "J6" = 291.0645
    "K6" = 442.8
    "L6" = 873.1935000000001
    "M6" = -330.8
    "N6" = -1097.1935
    "H18" = 797.4
    "I18" = 621.75
    "K18" = 621.75
    "M18" = -337.75
    "H97" = 1750.4615
    "J97" = 1623.3636
    "L97" = 4870.0908
    "N97" = -5862.0908
    "H129" = 3810.7
    "I129" = 707.46155
    "K129" = 2122.38465
    "M129" = 2877.61535
    "H132" = 1739.1177
    "I132" = 1752.86
    "K132" = 5258.58
    "M132" = -2728.58
    "H137" = 2592.1943
    "I137" = 3243.2222
    "K137" = 9729.6666
    "M137" = -7179.6666
    "H138" = 2611.34
    "I138" = 1042.5294
    "J138" = 2932.6626
    "K138" = 3127.5882
    "L138" = 8797.987800000001
    "M138" = 2012.4118
    "N138" = -19077.9878
}

$sheetUpdates["ARM"] = @{
    "H32" = 7592.55
    "I32" = 7224.726
    "K32" = 7224.726
    "M32" = -6937.726
    "H88" = 1747.7693
    "I88" = 5000
    "J88" = 1476.75
    "K88" = 5000
    "L88" = 1476.75
    "M88" = -4594
    "N88" = -2288.75
    "H91" = 1747.7693
    "I91" = 5000
    "J91" = 1476.75
    "K91" = 5000
    "L91" = 1476.75
    "M91" = -3596
    "N91" = -4284.75
    "H97" = 2415.4062
    "I97" = 868.1
    "J97" = 25625
    "K97" = 868.1
    "L97" = 25625
    "M97" = -372.1
    "N97" = -26617
    "H107" = 92099.5
    "J107" = 92099.5
    "L107" = 92099.5
    "N107" = -99779.5
    "H111" = 89994.5
    "J111" = 89994.5
    "L111" = 89994.5
    "N111" = -98174.5
    "H112" = 20025
    "J112" = 20025
    "L112" = 20025
    "N112" = -22979
    "H122" = 6986.5
    "I122" = 6579.8
    "K122" = 19739.4
    "M122" = -17289.4
    "H132" = 1609.5577
    "I132" = 1472.8049
    "J132" = 2119.2727
    "K132" = 4418.4147
    "L132" = 6357.8181
    "M132" = -1888.4147
    "N132" = -11417.8181
}

$sheetUpdates["BSM"] = @{
    "H20" = 1854.6809
    "I20" = 1738.5428
    "J20" = 2193.4167
    "K20" = 1738.5428
    "L20" = 2193.4167
    "M20" = -1491.5428
    "N20" = -2687.4167
}

$sheetUpdates["CRP"] = @{
    "H58" = 3141.96
    "I58" = 2627.8333
    "K58" = 2627.8333
    "M58" = -2424.8333
    "H102" = 65194.332
    "J102" = 65194.332
    "L102" = 65194.332
    "N102" = -70062.33199999999
    "H122" = 4314.5557
    "I122" = 4143.533
    "K122" = 12430.599
    "M122" = -9980.599000000002
    "H132" = 977946.5
    "I132" = 1740492.5
    "J132" = 3582.2222
    "K132" = 5221477.5
    "L132" = 10746.6666
    "M132" = -5218947.5
    "N132" = -15806.6666
    "H134" = 3438.6316
    "I134" = 959.5405
    "J134" = 8024.95
    "K134" = 2878.6215
    "L134" = 24074.85
    "M134" = -343.6214999999997
    "N134" = -29144.85
    "H136" = 3141.96
    "I136" = 2627.8333
    "K136" = 7883.499899999999
    "M136" = -5333.499899999999
}

$sheetUpdates["CUL"] = @{
    "H5" = 8919.857
    "I5" = 1797.6666
    "J5" = 10862.272
    "K5" = 5392.9998
    "L5" = 32586.816
    "M5" = -5280.9998
    "N5" = -32810.81600000001
    "H12" = 1626.0834
    "J12" = 977.93335
    "L12" = 2933.80005
    "N12" = -3279.80005
    "H23" = 71429460
    "J23" = 100001180
    "L23" = 300003540
    "N23" = -300004010
    "H135" = 8919.857
    "I135" = 1797.6666
    "J135" = 10862.272
    "K135" = 16178.9994
    "L135" = 97760.448
    "M135" = -13643.9994
    "N135" = -102830.448
}

$sheetUpdates["GSM"] = @{
    "H70" = 4263.636
    "I70" = 3732.318
    "J70" = 5326.273
    "K70" = 3732.318
    "L70" = 5326.273
    "M70" = -3462.318
    "N70" = -5866.273
    "H73" = 4263.636
    "I73" = 3732.318
    "J73" = 5326.273
    "K73" = 3732.318
    "L73" = 5326.273
    "M73" = -2796.318
    "N73" = -7198.273
    "H80" = 3248.0476
    "I80" = 2848.6667
    "K80" = 2848.6667
    "M80" = -1850.6667
    "H83" = 3248.0476
    "I83" = 2848.6667
    "K83" = 14243.3335
    "M83" = -9251.333500000001
    "H102" = 24824.783
    "J102" = 59018
    "L102" = 59018
    "N102" = -62262
    "H122" = 3000
    "J122" = 3000
    "L122" = 9000
    "N122" = -13900
    "H131" = 89500
    "J131" = 89500
    "L131" = 89500
    "N131" = -99580
    "H132" = 2952.55
    "I132" = 2964.4285
    "K132" = 8893.2855
    "M132" = -6363.2855
    "H136" = 39984.316
    "J136" = 39984.316
    "L136" = 119952.948
    "N136" = -125052.948
}

$sheetUpdates["LTW"] = @{
    "H7" = 5800.8335
    "I7" = 6144.6665
    "J7" = 4769.3335
    "K7" = 6144.6665
    "L7" = 4769.3335
    "M7" = -6032.6665
    "N7" = -4993.3335
    "H68" = 3586.5557
    "I68" = 3539.8572
    "K68" = 3539.8572
    "M68" = -2790.8572
    "H71" = 3586.5557
    "I71" = 3539.8572
    "K71" = 17699.286
    "M71" = -13955.286
    "H126" = 5800.8335
    "I126" = 6144.6665
    "J126" = 4769.3335
    "K126" = 18433.9995
    "L126" = 14308.0005
    "M126" = -15963.9995
    "N126" = -19248.0005
    "H136" = 3924.3713
    "I136" = 3329.4138
    "J136" = 6800
    "K136" = 9988.241399999999
    "L136" = 20400
    "M136" = -7438.241399999999
    "N136" = -25500
}

$sheetUpdates["WVR"] = @{
    "H81" = 22682.98
    "I81" = 43633.082
    "J81" = 4060.6667
    "K81" = 87266.164
    "L81" = 8121.3334
    "M81" = -86205.164
    "N81" = -10243.3334
    "H84" = 22682.98
    "I84" = 43633.082
    "J84" = 4060.6667
    "K84" = 436330.82
    "L84" = 40606.667
    "M84" = -431026.82
    "N84" = -51214.667
    "H122" = 3318.7646
    "I122" = 3185.4614
    "K122" = 9556.3842
    "M122" = -7106.3842
    "H132" = 1150.38
    "I132" = 1037.0769
    "K132" = 3111.2307
    "M132" = -581.2307000000001
    "H136" = 4065.2173
    "I136" = 4549.8706
    "J136" = 2320.4666
    "K136" = 13649.6118
    "L136" = 6961.399800000001
    "M136" = -11099.6118
    "N136" = -12061.3998
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $sheetUpdates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}

Write-Host "Updated $($sheetUpdates.Keys.Count) sheets."
